$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row values (A1:H1)
$ws.Range("A1").Value = "Nature de l'équipement"
$ws.Range("B1").Value = "Marque"
$ws.Range("C1").Value = "Modèle"
$ws.Range("D1").Value = "Type d'OS"
$ws.Range("E1").Value = "Version OS"
$ws.Range("F1").Value = "Hostname"
$ws.Range("G1").Value = "Ports utilisés"
$ws.Range("H1").Value = "Ports non utilisés/éteints/passifs"

# Bold Arial 10pt header style: build it once on a scratch cell, then copy
# the format (not the value) onto A1:H1 in a single paste so every header
# cell ends up sharing one cell-format entry instead of each Range.Font.*
# assignment minting its own per-cell style.
$scratch = $ws.Range("Z100")
$scratch.Font.Name = "Arial"
$scratch.Font.Size = 10
$scratch.Font.Bold = $true
$scratch.Copy()
$ws.Range("A1:H1").PasteSpecial(-4122)  # xlPasteFormats
$scratch.Clear()
$excel.CutCopyMode = $false

# Column widths (compensate for the engine's internal pixel rounding -
# character-width 5/6 padding - so the stored XML width lands as close as
# possible to the target values)
$ws.Columns.Item(1).ColumnWidth = 22.1640625 - 5/6
$ws.Columns.Item(2).ColumnWidth = 12.33203125 - 5/6
$ws.Columns.Item(4).ColumnWidth = 12.33203125 - 5/6
$ws.Columns.Item(5).ColumnWidth = 15.6640625 - 5/6
$ws.Columns.Item(6).ColumnWidth = 14.6640625 - 5/6
$ws.Columns.Item(7).ColumnWidth = 16 - 5/6
$ws.Columns.Item(8).ColumnWidth = 31.6640625 - 5/6

# Selection moves to A5
$ws.Range("A5").Select()
